$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1659.8
$ws.Range("I2").Value = 299.66666
$ws.Range("K2").Value = 299.66666
$ws.Range("M2").Value = -186.66666
$ws.Range("H6").Value = 785.1818
$ws.Range("I6").Value = 485.75
$ws.Range("K6").Value = 1457.25
$ws.Range("M6").Value = -1345.25
$ws.Range("H9").Value = 323.6875
$ws.Range("I9").Value = 390.41666
$ws.Range("J9").Value = 123.5
$ws.Range("K9").Value = 390.41666
$ws.Range("L9").Value = 123.5
$ws.Range("M9").Value = -221.41666
$ws.Range("N9").Value = -461.5
$ws.Range("H28").Value = 1309
$ws.Range("I28").Value = 755.4167
$ws.Range("J28").Value = 3523.3333
$ws.Range("K28").Value = 755.4167
$ws.Range("L28").Value = 3523.3333
$ws.Range("M28").Value = -270.4167
$ws.Range("N28").Value = -4493.3333
$ws.Range("H40").Value = 1784.1428
$ws.Range("I40").Value = 1784.1428
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1784.1428
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1609.1428
$ws.Range("H55").Value = 313.16666
$ws.Range("I55").Value = 199.66667
$ws.Range("K55").Value = 199.66667
$ws.Range("M55").Value = 14.33332999999999
$ws.Range("H62").Value = 7939.375
$ws.Range("I62").Value = 7846.625
$ws.Range("J62").Value = 8032.125
$ws.Range("K62").Value = 7846.625
$ws.Range("L62").Value = 8032.125
$ws.Range("M62").Value = -7222.625
$ws.Range("N62").Value = -9280.125
$ws.Range("H65").Value = 7939.375
$ws.Range("I65").Value = 7846.625
$ws.Range("J65").Value = 8032.125
$ws.Range("K65").Value = 39233.125
$ws.Range("L65").Value = 40160.625
$ws.Range("M65").Value = -36113.125
$ws.Range("N65").Value = -46400.625
$ws.Range("H74").Value = 2224.4443
$ws.Range("I74").Value = 2224.4443
$ws.Range("K74").Value = 2224.4443
$ws.Range("M74").Value = -1288.4443
$ws.Range("H77").Value = 2224.4443
$ws.Range("I77").Value = 2224.4443
$ws.Range("K77").Value = 11122.2215
$ws.Range("M77").Value = -6442.2215
$ws.Range("H81").Value = 43799
$ws.Range("J81").Value = 43799
$ws.Range("L81").Value = 43799
$ws.Range("N81").Value = -45795
$ws.Range("H84").Value = 43799
$ws.Range("J84").Value = 43799
$ws.Range("L84").Value = 131397
$ws.Range("N84").Value = -141381
$ws.Range("H87").Value = 137222.22
$ws.Range("J87").Value = 142500
$ws.Range("L87").Value = 142500
$ws.Range("N87").Value = -144996
$ws.Range("H88").Value = 5891.4165
$ws.Range("I88").Value = 884
$ws.Range("K88").Value = 884
$ws.Range("M88").Value = -478
$ws.Range("H90").Value = 137222.22
$ws.Range("J90").Value = 142500
$ws.Range("L90").Value = 427500
$ws.Range("N90").Value = -439980
$ws.Range("H91").Value = 5891.4165
$ws.Range("I91").Value = 884
$ws.Range("K91").Value = 884
$ws.Range("M91").Value = 520
$ws.Range("H92").Value = 1640.1904
$ws.Range("I92").Value = 303.8125
$ws.Range("J92").Value = 5916.6
$ws.Range("K92").Value = 303.8125
$ws.Range("L92").Value = 5916.6
$ws.Range("M92").Value = 944.1875
$ws.Range("N92").Value = -8412.6
$ws.Range("H94").Value = 2526.4167
$ws.Range("I94").Value = 2526.4167
$ws.Range("K94").Value = 2526.4167
$ws.Range("M94").Value = -2075.4167
$ws.Range("H96").Value = 4763285.5
$ws.Range("I96").Value = 7143928.5
$ws.Range("J96").Value = 1999.4
$ws.Range("K96").Value = 21431785.5
$ws.Range("L96").Value = 5998.200000000001
$ws.Range("M96").Value = -21430412.5
$ws.Range("N96").Value = -8744.200000000001
$ws.Range("H97").Value = 954.25
$ws.Range("J97").Value = 986.3333
$ws.Range("L97").Value = 2958.9999
$ws.Range("N97").Value = -3950.9999
$ws.Range("H99").Value = 547
$ws.Range("I99").Value = 295
$ws.Range("J99").Value = 727
$ws.Range("K99").Value = 885
$ws.Range("L99").Value = 2181
$ws.Range("M99").Value = 613
$ws.Range("N99").Value = -5177
$ws.Range("H100").Value = 42848.535
$ws.Range("I100").Value = 61203.2
$ws.Range("J100").Value = 6139.2
$ws.Range("K100").Value = 61203.2
$ws.Range("L100").Value = 6139.2
$ws.Range("M100").Value = -60662.2
$ws.Range("N100").Value = -7221.2
$ws.Range("H101").Value = 900
$ws.Range("I101").Value = 200
$ws.Range("J101").Value = 1075
$ws.Range("K101").Value = 600
$ws.Range("L101").Value = 3225
$ws.Range("M101").Value = 1022
$ws.Range("N101").Value = -6469
$ws.Range("H103").Value = 2999.5
$ws.Range("I103").Value = 999.6667
$ws.Range("J103").Value = 4999.3335
$ws.Range("K103").Value = 2999.0001
$ws.Range("L103").Value = 14998.0005
$ws.Range("M103").Value = -2413.0001
$ws.Range("N103").Value = -16170.0005
$ws.Range("H107").Value = 1303.1578
$ws.Range("J107").Value = 1300
$ws.Range("L107").Value = 1300
$ws.Range("N107").Value = -5140
$ws.Range("H125").Value = 76254.75
$ws.Range("I125").Value = 150265.5
$ws.Range("K125").Value = 1352389.5
$ws.Range("M125").Value = -1349929.5
$ws.Range("H135").Value = 16251.5
$ws.Range("I135").Value = 1539.8276
$ws.Range("J135").Value = 55036.816
$ws.Range("K135").Value = 13858.4484
$ws.Range("L135").Value = 495331.344
$ws.Range("M135").Value = -11323.4484
$ws.Range("N135").Value = -500401.344
$ws.Range("H137").Value = 10555.541
$ws.Range("I137").Value = 18543.277
$ws.Range("K137").Value = 55629.83099999999
$ws.Range("M137").Value = -53079.83099999999
$ws.Range("H141").Value = 2000
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 287533.28
$ws.Range("J6").Value = 287533.28
$ws.Range("L6").Value = 287533.28
$ws.Range("N6").Value = -287879.28
$ws.Range("H32").Value = 20360.695
$ws.Range("I32").Value = 21375.4
$ws.Range("J32").Value = 6408.5
$ws.Range("K32").Value = 21375.4
$ws.Range("L32").Value = 6408.5
$ws.Range("M32").Value = -21088.4
$ws.Range("N32").Value = -6982.5
$ws.Range("H37").Value = 29999.285
$ws.Range("I37").Value = 15000
$ws.Range("J37").Value = 49998.332
$ws.Range("K37").Value = 15000
$ws.Range("L37").Value = 49998.332
$ws.Range("M37").Value = -14727
$ws.Range("N37").Value = -50544.332
$ws.Range("H61").Value = 2254.2222
$ws.Range("I61").Value = 1725.3
$ws.Range("K61").Value = 1725.3
$ws.Range("M61").Value = -1513.3
$ws.Range("H74").Value = 353885.94
$ws.Range("I74").Value = 429461.72
$ws.Range("K74").Value = 429461.72
$ws.Range("M74").Value = -428587.72
$ws.Range("H77").Value = 353885.94
$ws.Range("I77").Value = 429461.72
$ws.Range("K77").Value = 2147308.6
$ws.Range("M77").Value = -2142940.6
$ws.Range("H102").Value = 3671.2285
$ws.Range("I102").Value = 2981.4814
$ws.Range("J102").Value = 5999.125
$ws.Range("K102").Value = 2981.4814
$ws.Range("L102").Value = 5999.125
$ws.Range("M102").Value = -1359.4814
$ws.Range("N102").Value = -9243.125
$ws.Range("H110").Value = 2849.5
$ws.Range("I110").Value = 2849.5
$ws.Range("K110").Value = 2849.5
$ws.Range("M110").Value = -804.5
$ws.Range("H122").Value = 2500.4473
$ws.Range("I122").Value = 2389.361
$ws.Range("J122").Value = 4500
$ws.Range("K122").Value = 7168.083
$ws.Range("L122").Value = 13500
$ws.Range("M122").Value = -4718.083
$ws.Range("N122").Value = -18400
$ws.Range("H132").Value = 1604.0322
$ws.Range("I132").Value = 1081.2693
$ws.Range("K132").Value = 3243.8079
$ws.Range("M132").Value = -713.8078999999998
$ws.Range("H136").Value = 2254.2222
$ws.Range("I136").Value = 1725.3
$ws.Range("K136").Value = 5175.9
$ws.Range("M136").Value = -2625.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 73433.28999999999
$ws.Range("I20").Value = 85422.164
$ws.Range("J20").Value = 1500
$ws.Range("K20").Value = 85422.164
$ws.Range("L20").Value = 1500
$ws.Range("M20").Value = -85175.164
$ws.Range("N20").Value = -1994
$ws.Range("H22").Value = 348.66666
$ws.Range("I22").Value = 318.4
$ws.Range("K22").Value = 318.4
$ws.Range("M22").Value = -145.4
$ws.Range("H86").Value = 1802.9375
$ws.Range("I86").Value = 1523.2
$ws.Range("J86").Value = 5999
$ws.Range("K86").Value = 1523.2
$ws.Range("L86").Value = 5999
$ws.Range("M86").Value = -400.2
$ws.Range("N86").Value = -8245
$ws.Range("H89").Value = 1802.9375
$ws.Range("I89").Value = 1523.2
$ws.Range("J89").Value = 5999
$ws.Range("K89").Value = 7616
$ws.Range("L89").Value = 29995
$ws.Range("M89").Value = -2000
$ws.Range("N89").Value = -41227
$ws.Range("H92").Value = 1950
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("H94").Value = 2133
$ws.Range("I94").Value = 2133
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 2133
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1682
$ws.Range("H97").Value = 16500
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("H99").Value = 1600
$ws.Range("I99").Value = 1120.6154
$ws.Range("J99").Value = 2292.4443
$ws.Range("K99").Value = 1120.6154
$ws.Range("L99").Value = 2292.4443
$ws.Range("M99").Value = 377.3846000000001
$ws.Range("N99").Value = -5288.4443
$ws.Range("H105").Value = 1544.091
$ws.Range("I105").Value = 1057.8334
$ws.Range("J105").Value = 3732.25
$ws.Range("K105").Value = 1057.8334
$ws.Range("L105").Value = 3732.25
$ws.Range("M105").Value = 689.1666
$ws.Range("N105").Value = -7226.25
$ws.Range("H107").Value = 31834.666
$ws.Range("I107").Value = 37281.668
$ws.Range("K107").Value = 37281.668
$ws.Range("M107").Value = -35361.668
$ws.Range("N92").ClearContents()
$ws.Range("N94").ClearContents()
$ws.Range("M97").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 44325
$ws.Range("J14").Value = 44325
$ws.Range("L14").Value = 44325
$ws.Range("N14").Value = -44665
$ws.Range("H16").Value = 1259.2727
$ws.Range("I16").Value = 984.8421
$ws.Range("J16").Value = 2997.3333
$ws.Range("K16").Value = 984.8421
$ws.Range("L16").Value = 2997.3333
$ws.Range("M16").Value = -697.8421
$ws.Range("N16").Value = -3571.3333
$ws.Range("H31").Value = 2274841.5
$ws.Range("I31").Value = 2779805.5
$ws.Range("J31").Value = 2503.625
$ws.Range("K31").Value = 2779805.5
$ws.Range("L31").Value = 2503.625
$ws.Range("M31").Value = -2779510.5
$ws.Range("N31").Value = -3093.625
$ws.Range("H34").Value = 2274841.5
$ws.Range("I34").Value = 2779805.5
$ws.Range("J34").Value = 2503.625
$ws.Range("K34").Value = 2779805.5
$ws.Range("L34").Value = 2503.625
$ws.Range("M34").Value = -2779603.5
$ws.Range("N34").Value = -2907.625
$ws.Range("H51").Value = 30924.5
$ws.Range("J51").Value = 30924.5
$ws.Range("L51").Value = 30924.5
$ws.Range("N51").Value = -32396.5
$ws.Range("H61").Value = 30924.5
$ws.Range("J61").Value = 30924.5
$ws.Range("L61").Value = 30924.5
$ws.Range("N61").Value = -31620.5
$ws.Range("H70").Value = 46250
$ws.Range("J70").Value = 46250
$ws.Range("L70").Value = 46250
$ws.Range("N70").Value = -46880
$ws.Range("H73").Value = 46250
$ws.Range("J73").Value = 46250
$ws.Range("L73").Value = 46250
$ws.Range("N73").Value = -48434
$ws.Range("H107").Value = 653.5714
$ws.Range("I107").Value = 483.5263
$ws.Range("J107").Value = 1012.55554
$ws.Range("K107").Value = 483.5263
$ws.Range("L107").Value = 1012.55554
$ws.Range("M107").Value = 1436.4737
$ws.Range("N107").Value = -4852.55554
$ws.Range("H113").Value = 1259.2727
$ws.Range("I113").Value = 984.8421
$ws.Range("J113").Value = 2997.3333
$ws.Range("K113").Value = 984.8421
$ws.Range("L113").Value = 2997.3333
$ws.Range("M113").Value = 1185.1579
$ws.Range("N113").Value = -7337.3333
$ws.Range("H122").Value = 2323.4285
$ws.Range("I122").Value = 2212.75
$ws.Range("J122").Value = 2471
$ws.Range("K122").Value = 6638.25
$ws.Range("L122").Value = 7413
$ws.Range("M122").Value = -4188.25
$ws.Range("N122").Value = -12313
$ws.Range("H132").Value = 25541.48
$ws.Range("I132").Value = 28832.273
$ws.Range("J132").Value = 1409
$ws.Range("K132").Value = 86496.819
$ws.Range("L132").Value = 4227
$ws.Range("M132").Value = -83966.819
$ws.Range("N132").Value = -9287
$ws.Range("H134").Value = 2279.25
$ws.Range("I134").Value = 2223.8076
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 6671.4228
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -4136.4228
$ws.Range("N134").Value = -14070
$ws.Range("H138").Value = 112301.5
$ws.Range("J138").Value = 112301.5
$ws.Range("L138").Value = 112301.5
$ws.Range("N138").Value = -122581.5
$ws.Range("H141").Value = 296391.78
$ws.Range("J141").Value = 296391.78
$ws.Range("L141").Value = 296391.78
$ws.Range("N141").Value = -306751.78

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 771710.5600000001
$ws.Range("I2").Value = 868162.4399999999
$ws.Range("J2").Value = 95.5
$ws.Range("K2").Value = 5208974.64
$ws.Range("L2").Value = 573
$ws.Range("M2").Value = -5208861.64
$ws.Range("N2").Value = -799
$ws.Range("H94").Value = 1999.5
$ws.Range("I94").Value = 1999.5
$ws.Range("K94").Value = 5998.5
$ws.Range("M94").Value = -5322.5
$ws.Range("H109").Value = 3106.05
$ws.Range("I109").Value = 1990.75
$ws.Range("J109").Value = 3849.5833
$ws.Range("K109").Value = 5972.25
$ws.Range("L109").Value = 11548.7499
$ws.Range("M109").Value = -4932.25
$ws.Range("N109").Value = -13628.7499
$ws.Range("H120").Value = 5000
$ws.Range("I120").Value = 5000
$ws.Range("K120").Value = 15000
$ws.Range("M120").Value = -10162
$ws.Range("H129").Value = 4103.7144
$ws.Range("I129").Value = 1917.6
$ws.Range("K129").Value = 5752.799999999999
$ws.Range("M129").Value = -752.7999999999993
$ws.Range("H131").Value = 144106.33
$ws.Range("J131").Value = 2610.3076
$ws.Range("L131").Value = 7830.9228
$ws.Range("N131").Value = -17910.9228
$ws.Range("H134").Value = 1622.75
$ws.Range("I134").Value = 1622.75
$ws.Range("K134").Value = 4868.25
$ws.Range("M134").Value = 201.75
$ws.Range("H138").Value = 9750.166999999999
$ws.Range("J138").Value = 9141.857
$ws.Range("L138").Value = 27425.571
$ws.Range("N138").Value = -37705.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 297.14285
$ws.Range("I2").Value = 253.8
$ws.Range("K2").Value = 253.8
$ws.Range("M2").Value = -140.8
$ws.Range("H4").Value = 5000
$ws.Range("J4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5224
$ws.Range("H36").Value = 10833.333
$ws.Range("I36").Value = 1000
$ws.Range("J36").Value = 15750
$ws.Range("K36").Value = 1000
$ws.Range("L36").Value = 15750
$ws.Range("M36").Value = -515
$ws.Range("N36").Value = -16720
$ws.Range("H80").Value = 6993.222
$ws.Range("I80").Value = 4237.8
$ws.Range("J80").Value = 10437.5
$ws.Range("K80").Value = 4237.8
$ws.Range("L80").Value = 10437.5
$ws.Range("M80").Value = -3239.8
$ws.Range("N80").Value = -12433.5
$ws.Range("H83").Value = 6993.222
$ws.Range("I83").Value = 4237.8
$ws.Range("J83").Value = 10437.5
$ws.Range("K83").Value = 21189
$ws.Range("L83").Value = 52187.5
$ws.Range("M83").Value = -16197
$ws.Range("N83").Value = -62171.5
$ws.Range("H97").Value = 1618.3043
$ws.Range("I97").Value = 1339.3158
$ws.Range("J97").Value = 2943.5
$ws.Range("K97").Value = 1339.3158
$ws.Range("L97").Value = 2943.5
$ws.Range("M97").Value = -843.3158000000001
$ws.Range("N97").Value = -3935.5
$ws.Range("H102").Value = 22503.28
$ws.Range("I102").Value = 25022.227
$ws.Range("K102").Value = 25022.227
$ws.Range("M102").Value = -23400.227
$ws.Range("H122").Value = 2982.1052
$ws.Range("I122").Value = 2758.889
$ws.Range("J122").Value = 7000
$ws.Range("K122").Value = 8276.667000000001
$ws.Range("L122").Value = 21000
$ws.Range("M122").Value = -5826.667000000001
$ws.Range("N122").Value = -25900
$ws.Range("H126").Value = 2136.6
$ws.Range("I126").Value = 1978.6154
$ws.Range("J126").Value = 2430
$ws.Range("K126").Value = 5935.8462
$ws.Range("L126").Value = 7290
$ws.Range("M126").Value = -3465.8462
$ws.Range("N126").Value = -12230
$ws.Range("H132").Value = 5249.552
$ws.Range("I132").Value = 5902.64
$ws.Range("J132").Value = 1167.75
$ws.Range("K132").Value = 17707.92
$ws.Range("L132").Value = 3503.25
$ws.Range("M132").Value = -15177.92
$ws.Range("N132").Value = -8563.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1423.6111
$ws.Range("I22").Value = 892.5714
$ws.Range("J22").Value = 1761.5454
$ws.Range("K22").Value = 892.5714
$ws.Range("L22").Value = 1761.5454
$ws.Range("M22").Value = -597.5714
$ws.Range("N22").Value = -2351.5454
$ws.Range("H27").Value = 1423.6111
$ws.Range("I27").Value = 892.5714
$ws.Range("J27").Value = 1761.5454
$ws.Range("K27").Value = 892.5714
$ws.Range("L27").Value = 1761.5454
$ws.Range("M27").Value = -785.5714
$ws.Range("N27").Value = -1975.5454
$ws.Range("H46").Value = 3598.8572
$ws.Range("I46").Value = 792.4
$ws.Range("J46").Value = 6837.077
$ws.Range("K46").Value = 792.4
$ws.Range("L46").Value = 6837.077
$ws.Range("M46").Value = -604.4
$ws.Range("N46").Value = -7213.077
$ws.Range("H55").Value = 978.9259
$ws.Range("I55").Value = 362.9091
$ws.Range("J55").Value = 1402.4375
$ws.Range("K55").Value = 362.9091
$ws.Range("L55").Value = 1402.4375
$ws.Range("M55").Value = -189.9091
$ws.Range("N55").Value = -1748.4375
$ws.Range("H56").Value = 54003.2
$ws.Range("J56").Value = 54003.2
$ws.Range("L56").Value = 54003.2
$ws.Range("N56").Value = -55385.2
$ws.Range("H74").Value = 56999.285
$ws.Range("I74").Value = 35000
$ws.Range("J74").Value = 65799
$ws.Range("K74").Value = 35000
$ws.Range("L74").Value = 65799
$ws.Range("M74").Value = -34002
$ws.Range("N74").Value = -67795
$ws.Range("H77").Value = 56999.285
$ws.Range("I77").Value = 35000
$ws.Range("J77").Value = 65799
$ws.Range("K77").Value = 105000
$ws.Range("L77").Value = 197397
$ws.Range("M77").Value = -100008
$ws.Range("N77").Value = -207381
$ws.Range("H122").Value = 7446.9536
$ws.Range("I122").Value = 7560.9
$ws.Range("J122").Value = 7184
$ws.Range("K122").Value = 22682.7
$ws.Range("L122").Value = 21552
$ws.Range("M122").Value = -20232.7
$ws.Range("N122").Value = -26452
$ws.Range("H136").Value = 9382.591
$ws.Range("I136").Value = 9395.85
$ws.Range("K136").Value = 28187.55
$ws.Range("M136").Value = -25637.55

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 29969
$ws.Range("J15").Value = 29969
$ws.Range("L15").Value = 29969
$ws.Range("N15").Value = -30545
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("H62").Value = 7330.6665
$ws.Range("I62").Value = 6374.75
$ws.Range("K62").Value = 6374.75
$ws.Range("M62").Value = -5750.75
$ws.Range("H65").Value = 7330.6665
$ws.Range("I65").Value = 6374.75
$ws.Range("K65").Value = 31873.75
$ws.Range("M65").Value = -28753.75
$ws.Range("H75").Value = 105000
$ws.Range("J75").Value = 105000
$ws.Range("L75").Value = 105000
$ws.Range("N75").Value = -106872
$ws.Range("H78").Value = 105000
$ws.Range("J78").Value = 105000
$ws.Range("L78").Value = 315000
$ws.Range("N78").Value = -324360
$ws.Range("H81").Value = 3976
$ws.Range("I81").Value = 4461.4165
$ws.Range("K81").Value = 8922.833000000001
$ws.Range("M81").Value = -7861.833000000001
$ws.Range("H84").Value = 3976
$ws.Range("I84").Value = 4461.4165
$ws.Range("K84").Value = 44614.165
$ws.Range("M84").Value = -39310.165
$ws.Range("H96").Value = 1827.5555
$ws.Range("I96").Value = 1346.2
$ws.Range("J96").Value = 2429.25
$ws.Range("K96").Value = 1346.2
$ws.Range("L96").Value = 2429.25
$ws.Range("M96").Value = 26.79999999999995
$ws.Range("N96").Value = -5175.25
$ws.Range("H126").Value = 360370.16
$ws.Range("I126").Value = 3431.8333
$ws.Range("K126").Value = 10295.4999
$ws.Range("M126").Value = -7825.499899999999
$ws.Range("H132").Value = 3535.3333
$ws.Range("I132").Value = 3492.138
$ws.Range("K132").Value = 10476.414
$ws.Range("M132").Value = -7946.414000000001
$ws.Range("H136").Value = 25742.777
$ws.Range("I136").Value = 27642.24
$ws.Range("J136").Value = 1999.5
$ws.Range("K136").Value = 82926.72
$ws.Range("L136").Value = 5998.5
$ws.Range("M136").Value = -80376.72
$ws.Range("N136").Value = -11098.5
$ws.Range("M61").ClearContents()
